# Applies the LOQ4002 content corrections:
#  - Gives "Objetivos:" its own proper text (was incorrectly showing the docente's name)
#  - Inserts a new row (13) holding the "Docentes responsaveis" data, shifting the
#    remaining rows down by one
#  - Fixes each of the now-shifted rows so the data in column B/C matches the label
#    in column A again (Programa resumido, Programa, Metodo, Criterio, Norma de
#    recuperacao, Bibliografia all get their correct text)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text constants -------------------------------------------------------

$OBJETIVOS = "Capacitar os alunos a calcular os parâmetros de projeto de reatores ideais, a distinguir entre um reator ideal e um real, e a compreender a influência da temperatura e pressão no projeto de reatores químicos."

$DOCENTE = "5963230 - Leandro Gonçalves de Aguiar"

$PROGRAMA_RESUMIDO = "1. Introdução a Reatores. 2. Modelos Ideais de Reatores Químicos Isotérmicos  Reações Simples. 3. Reações Múltiplas em Reatores Ideais. 4. Efeitos Térmicos em Reatores Ideais. 5. Reatores Catalíticos Heterogêneos. 6. Reatores Não-Ideais"

$PROGRAMA = "1. Introdução a Reatores: Conceitos básicos`n2. Modelos Ideais de Reatores Químicos Isotérmicos  Reações Simples: `n2.1) Equações fundamentais de projeto de reatores`n2.2) Reator tanque descontínuo (BSTR)`n2.3) Reator tanque de mistura contínuo (CSTR)`n2.4) Reator tubular de fluxo pistonado (PFR)`n2.5) Comparação de desempenho de reatores CSTR e PFR`n2.6) Reatores CSTR em cascata`n2.7) Associação mista de reatores em série: CSTR e PFR`n2.8) Reatores com reciclo`n2.9) Reações auto-catalíticas`n2.10) Reatores semi-contínuos`n3. Reações Múltiplas em Reatores Ideais`n3.1) Noções gerais: otimização, rendimento e seletividade`n3.2) Reações paralelas e reações em série`n3.3) Sistemas com reações série-paralelo: reações de múltipla substituição e reações poliméricas`n3.4) Problemas simples de otimização`n4. Efeitos Térmicos em Reatores Ideais`n4.1) Equação do balanço de energia`n4.2) Balanço de energia aplicado ao BSTR`n4.3) Balanço de energia aplicado ao CSTR`n4.4) Balanço de energia aplicado ao PFR`n5. Reatores Catalíticos Heterogêneos`n5.1) Introdução`n5.2) Efeito dos processos físicos sobre a taxa de reação`n5.2.1  Fenômenos interfases`n5.2.2  Fenômenos intrapartícula`n5.2.3  Difusão e reação em catalisadores porosos`n5.3) Cálculo de reatores de leito fixo`n5.4) Reatores trifásicos`n6. Reatores Não-Ideais`n6.1) A distribuição dos tempos de residência`n6.2) Modelos dos tanques contínuos em série`n6.3) Modelo da dispersão axial"

$METODO = "Duas provas escritas e eventual apresentação de trabalho."

$CRITERIO = "Nota(N) = 50% Prova P1 + 50% Prova P2. Os pesos poderão ser redefinidos caso seja incorporada nota de trabalho."

$NORMA = "Média Final = (N + Prova Recuperação)/2"

$BIBLIOGRAFIA = "FOGLER, H. S. Elementos de Engenharia das Reações Químicas. 3. ed. Rio de Janeiro: LTC Editora, 2002.`nLEVENSPIEL, O. Chemical Reaction Engineering. 3. ed. New York: John Wiley & Sons, 1998.`nHILL, C.G. An Introduction to Chemical Engineering Kinetics and Reactor Design. New York: John Wiley&Sons, 1977.`nSMITH, J.M. Chemical Engineering Kinetics. 3rd. ed. New York :  McGraw-Hill, 1981.`nDENBIGH, K.; TURNER, R. Introduction to Chemical Reaction Design. Cambridge: Cambridge University Press, 1970.`nFROMENT, G.F.; BISCHOFF, K.B. Chemical Reactor Analysis And Design. 2nd ed.  New York: John Wiley & Sons, 1990."

# --- 1. Fix Objetivos text (row 10) ---------------------------------------

$ws.Range("B10:C10").Value = $OBJETIVOS

# --- 2. Insert a new row at 13 for the "Docentes responsaveis" data -------
#        only columns B:C are shifted down, so column A (which already has
#        the right labels for rows 13 onward) stays put.

$ws.Range("B13:C13").Insert(-4121)
$ws.Range("A13").Style = "Normal"

# The newly-materialised B13:C13 cells don't automatically pick up the
# normal "wrap text" column formatting, so copy it over explicitly from
# the (still correctly-formatted, untouched) row below.
$ws.Range("B14:C14").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B13:C13").Value = $DOCENTE

# --- 3. Re-fill the shifted rows with the matching text -------------------

$ws.Range("B14:C14").Value = $PROGRAMA_RESUMIDO
$ws.Range("B16:C16").Value = $PROGRAMA
$ws.Range("B19:C19").Value = $METODO
$ws.Range("B20:C20").Value = $CRITERIO
$ws.Range("B21:C21").Value = $NORMA
$ws.Range("B22:C22").Value = $BIBLIOGRAFIA

# Re-assigning long, wrapped, multi-line text can make Excel auto-fit the
# row height; pin the affected row back to its original fixed height.
$ws.Rows("16").RowHeight = 120
